# Update cryptocurrency price/volume snapshot data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.633.44"
$ws.Range("E2").Value = "  +0.99%  "
$ws.Range("D3").Value = "1.630.99"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.88"
$ws.Range("E5").Value = "  +0.37%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.499"
$ws.Range("E6").Value = "  +2.85%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +1.55%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0622"
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.19"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +3.42%  "
$ws.Range("D12").Value = "1.858.08"
$ws.Range("E12").Value = "  +0.54%  "
$ws.Range("D13").Value = "1.611.53"
$ws.Range("E13").Value = "  -1.04%  "
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.524"
$ws.Range("E15").Value = "  +0.99%  "
$ws.Range("D16").Value = "26.608.97"
$ws.Range("E16").Value = "  +0.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.28"
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "0.0₃0741"
$ws.Range("E18").Value = "  +1.90%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "218.36"
$ws.Range("E19").Value = "  +7.70%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.17"
$ws.Range("E22").Value = "  +2.16%  "
$ws.Range("E23").Value = "  +0.17%  "
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("E25").Value = "  +2.23%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  +0.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.84"
$ws.Range("E28").Value = "  +3.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.54"
$ws.Range("E29").Value = "  +2.24%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0504"
$ws.Range("E30").Value = "  -3.77%  "
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  +3.80%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +2.56%  "
$ws.Range("E34").Value = "  +0.32%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").Value = "1.213.28"
$ws.Range("E36").Value = "  +3.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0172"
$ws.Range("E37").Value = "  +4.74%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("E40").Value = "  +0.65%  "
$ws.Range("E41").Value = "  -2.06%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.791"
$ws.Range("E43").Value = "  +0.96%  "
$ws.Range("D44").Value = "1.776.15"
$ws.Range("E44").Value = "  +0.92%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "92.88"
$ws.Range("E45").Value = "  +0.19%  "
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "54.79"
$ws.Range("E47").Value = "  +1.76%  "
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0511"
$ws.Range("E48").Value = "  +0.52%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.59"
$ws.Range("E49").Value = "  +3.96%  "
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.409"
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.00"
$ws.Range("E51").Value = "  +0.22%  "
